# Pioneer Gliders Calibration and ingest CSV
# Update the glider Ref Des from CP05MOAS-GL003 to CP05MOAS-GL375 across the
# workbook (Moorings!A2 and every "CP05MOAS-GL003-xx-xxxxxxxxx" Ref Des on
# the Asset_Cal_Info sheet), then leave the same cells selected that the
# author's Excel session ended up on.

$wb = $excel.ActiveWorkbook

# Replace every occurrence of the old glider id with the new one, on every
# sheet, so both the standalone "CP05MOAS-GL003" cell and the longer
# "CP05MOAS-GL003-##-XXXXXXXXX" Ref Des strings get updated consistently.
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("GL003", "GL375")
}

# Moorings sheet: restore tab selection / active cell.
$wsMoorings = $wb.Worksheets.Item("Moorings")
$wsMoorings.Activate()
$wsMoorings.Range("D25").Select()

# Asset_Cal_Info sheet: restore its own active cell.
$wsAsset = $wb.Worksheets.Item("Asset_Cal_Info")
$wsAsset.Activate()
$wsAsset.Range("D19").Select()

# Leave the Moorings tab as the selected/active sheet, matching the source.
$wsMoorings.Activate()
